# Change the table style on the three tables (slides 14, 15, 16) from
# the default "Table_0" custom style ({7D36701E-6A0A-475B-A822-E8C1F818B1B3})
# to the built-in "Medium Style 2 - Accent 1" style
# ({015F5669-C73A-4331-A0DE-6870BCF93A80}).

$p = $ppt.ActivePresentation
$newStyleId = "{015F5669-C73A-4331-A0DE-6870BCF93A80}"
$targetSlides = @(14, 15, 16)

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
